$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 1.88
$ws.Range("R2").Value = 1.98
$ws.Range("S2").Value = 2.5
$ws.Range("T2").Value = 1.5
$ws.Range("U2").Value = 4
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 1.17
$ws.Range("G3").Value = 2.63
$ws.Range("I3").Value = 3.7
$ws.Range("N3").Value = 3.6
$ws.Range("AD3").Value = 10
$ws.Range("AM3").Value = 6
$ws.Range("AQ3").Value = 51
$ws.Range("G4").Value = 4.5
$ws.Range("I4").Value = 2
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("AC4").Value = 8.5
$ws.Range("AD4").Value = 21
$ws.Range("AE4").Value = 17
$ws.Range("AI4").Value = 5.5
$ws.Range("AJ4").Value = 6.5
$ws.Range("AK4").Value = 23
$ws.Range("AN4").Value = 7.5
$ws.Range("G5").Value = 1.85
$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 4.75
$ws.Range("J5").Value = 2.6
$ws.Range("L5").Value = 5.5
$ws.Range("Q5").Value = 1.93
$ws.Range("R5").Value = 1.93
$ws.Range("AA5").Value = 2.25
$ws.Range("AB5").Value = 1.57
$ws.Range("AC5").Value = 5
$ws.Range("AM5").Value = 9.5
$ws.Range("AO5").Value = 17
$ws.Range("G6").Value = 2.55
$ws.Range("I6").Value = 3.3
$ws.Range("J6").Value = 3.5
$ws.Range("L6").Value = 4.33
$ws.Range("R6").Value = 1.63
$ws.Range("AC6").Value = 5.5
$ws.Range("AD6").Value = 10
$ws.Range("AE6").Value = 11
$ws.Range("AF6").Value = 26
$ws.Range("AH6").Value = 41
$ws.Range("AL6").Value = 101
$ws.Range("AN6").Value = 15
$ws.Range("AP6").Value = 41
$ws.Range("U7").Value = 3.6
$ws.Range("V7").Value = 1.28
$ws.Range("M9").Value = 1.1
$ws.Range("N9").Value = 7
$ws.Range("Q9").Value = 1.93
$ws.Range("R9").Value = 1.93
$ws.Range("S9").Value = 2.5
$ws.Range("T9").Value = 1.5
$ws.Range("U9").Value = 4.1
$ws.Range("V9").Value = 1.22
$ws.Range("W9").Value = 5
$ws.Range("X9").Value = 1.17
$ws.Range("G11").Value = 3.35
$ws.Range("H11").Value = 2.6
$ws.Range("I11").Value = 2.52
$ws.Range("J11").Value = 4.15
$ws.Range("K11").Value = 1.75
$ws.Range("L11").Value = 3.35
$ws.Range("M11").Value = 1.17
$ws.Range("N11").Value = 4.4
$ws.Range("O11").Value = 1.7
$ws.Range("S11").Value = 3.05
$ws.Range("T11").Value = 1.33
$ws.Range("W11").Value = 5.6
$ws.Range("X11").Value = 1.11
$ws.Range("Z11").Value = 2.05
$ws.Range("AC11").Value = 6.4
$ws.Range("AD11").Value = 15
$ws.Range("AE11").Value = 13
$ws.Range("AF11").Value = 50
$ws.Range("AG11").Value = 45
$ws.Range("AH11").Value = 70
$ws.Range("AI11").Value = 4.4
$ws.Range("AJ11").Value = 5.5
$ws.Range("AM11").Value = 5.4
$ws.Range("AN11").Value = 10.5
$ws.Range("AO11").Value = 10.75
$ws.Range("AP11").Value = 29
$ws.Range("AQ11").Value = 30
$ws.Range("AR11").Value = 60
$ws.Range("G15").Value = 4.35
$ws.Range("H15").Value = 3.9
$ws.Range("J15").Value = 4.45
$ws.Range("K15").Value = 2.3
$ws.Range("L15").Value = 2.18
$ws.Range("AA15").Value = 1.62
$ws.Range("AB15").Value = 2.02
$ws.Range("AD15").Value = 27
$ws.Range("AE15").Value = 14
$ws.Range("AH15").Value = 37
$ws.Range("AJ15").Value = 7.8
$ws.Range("AK15").Value = 14
$ws.Range("AL15").Value = 55
$ws.Range("AM15").Value = 8.5
$ws.Range("AN15").Value = 8.75
$ws.Range("AS15").Value = 350
$ws.Range("G18").Value = 2.6
$ws.Range("H18").Value = 3.5
$ws.Range("I18").Value = 2.55
$ws.Range("J18").Value = 3
$ws.Range("K18").Value = 2.25
$ws.Range("L18").Value = 3
$ws.Range("M18").Value = 1.03
$ws.Range("N18").Value = 10
$ws.Range("S18").Value = 1.65
$ws.Range("T18").Value = 2.2
$ws.Range("W18").Value = 2.5
$ws.Range("X18").Value = 1.5
$ws.Range("AA18").Value = 1.53
$ws.Range("AB18").Value = 2.38
$ws.Range("AD18").Value = 15
$ws.Range("AE18").Value = 10
$ws.Range("AF18").Value = 26
$ws.Range("AG18").Value = 19
$ws.Range("AH18").Value = 23
$ws.Range("AI18").Value = 15
$ws.Range("AJ18").Value = 7
$ws.Range("AL18").Value = 34
$ws.Range("AM18").Value = 11
$ws.Range("AN18").Value = 15
$ws.Range("AP18").Value = 26
$ws.Range("AQ18").Value = 19
$ws.Range("G22").Value = 2.7
$ws.Range("H22").Value = 3.25
$ws.Range("I22").Value = 2.5
$ws.Range("J22").Value = 3.25
$ws.Range("L22").Value = 3.1
$ws.Range("N22").Value = 13
$ws.Range("AC22").Value = 10
$ws.Range("AE22").Value = 10
$ws.Range("AF22").Value = 26
$ws.Range("AP22").Value = 26
$ws.Range("G24").Value = 1.4
$ws.Range("H24").Value = 4.8
$ws.Range("I24").Value = 6.4
$ws.Range("J24").Value = 1.82
$ws.Range("K24").Value = 2.65
$ws.Range("L24").Value = 5.5
$ws.Range("P24").Value = 5.7
$ws.Range("T24").Value = 2.95
$ws.Range("X24").Value = 1.87
$ws.Range("AA24").Value = 1.5
$ws.Range("AB24").Value = 2.4
$ws.Range("AD24").Value = 9.75
$ws.Range("AF24").Value = 11.5
$ws.Range("AG24").Value = 10
$ws.Range("AH24").Value = 17
$ws.Range("AJ24").Value = 10.75
$ws.Range("AK24").Value = 14.5
$ws.Range("AL24").Value = 40
$ws.Range("AM24").Value = 29
$ws.Range("AN24").Value = 55
$ws.Range("AO24").Value = 20
$ws.Range("AP24").Value = 150
$ws.Range("AQ24").Value = 55
$ws.Range("AR24").Value = 40
$ws.Range("AS24").Value = 200
